$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.081.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.758.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.03"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.94"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.30%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.62%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.68"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -16.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.248.65"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.93"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.688.62"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.762.57"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.14"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "357.46"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.535"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.58"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.96%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.57"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = ("0.0{0}0910" -f $sub3)
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.32"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.97"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.21%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.08"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.26"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.01"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "345.25"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.35"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.20"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.20"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.45"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.85"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.33%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.632"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0255"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.04"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.18%  "
